$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: B27 changes from text "3" to a real number 3
$ws.Range("B27").Value = 3

# New row 28: copy of row27-like annotation but with new data
$ws.Range("A28").Value = "Ying Tang"
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "3"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "无"
$ws.Range("D28").Value = "FBK"
$ws.Range("E28").Value = "MET"
$ws.Range("F28").Value = "7e3af7f2-ed39-457b-b159-8a754cc477a9"
$ws.Range("G28").Value = "HkwZSG-CZ_annotated.xlsx"
$ws.Range("H28").Value = "In general, computational wall time of MoS is actually sub-linear w.r.t. the number of mixture components."
